$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.962.71"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "1.744.94"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.48"
$ws.Range("E5").Value = "  +5.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5044"
$ws.Range("E7").Value = "  -4.67%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2742"
$ws.Range("E8").Value = "  -2.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06186"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07273"
$ws.Range("D11").Value = "1.742.59"
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.6531"
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.15"
$ws.Range("E13").Value = "  -1.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.646"
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.67"
$ws.Range("E15").Value = "  -1.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9996"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9998"
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").Value = "25.974.45"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("E19").Value = "  +0.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006826"
$ws.Range("E20").Value = "  +0.84%  "
$ws.Range("D21").Value = "1.966.25"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.409"
$ws.Range("E22").Value = "  +1.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.708"
$ws.Range("E23").Value = "  -0.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.398"
$ws.Range("E24").Value = "  +3.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.72"
$ws.Range("E25").Value = "  -2.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.503"
$ws.Range("E26").Value = "  -1.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.24"
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.774"
$ws.Range("E28").Value = "  -1.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "105.56"
$ws.Range("E29").Value = "  +0.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.868"
$ws.Range("E30").Value = "  +2.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08215"
$ws.Range("E31").Value = "  -1.11%  "
$ws.Range("E32").Value = "  -0.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04686"
$ws.Range("E33").Value = "  +1.01%  "
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9928"
$ws.Range("E35").Value = "  -1.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6192"
$ws.Range("E36").Value = "  -2.10%  "
$ws.Range("E37").Value = "  +1.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01606"
$ws.Range("E38").Value = "  -1.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.923"
$ws.Range("E39").Value = "  -2.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9996"
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "100.40"
$ws.Range("E41").Value = "  -1.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.3934"
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7575"
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.015"
$ws.Range("E44").Value = "  -1.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1146"
$ws.Range("E45").Value = "  -0.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.301"
$ws.Range("E46").Value = "  -0.70%  "
$ws.Range("E47").Value = "  +1.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05276"
$ws.Range("E48").Value = "  -1.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.62"
$ws.Range("E49").Value = "  -1.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.555"
$ws.Range("E50").Value = "  -0.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3439"
$ws.Range("E51").Value = "  -1.05%  "
